$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.478.12'
$ws.Range("E2").Value = '  -1.08%  '
$ws.Range("D3").Value = '2.108.72'
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.89'
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5261'
$ws.Range("E7").Value = '  -1.32%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4519'
$ws.Range("E8").Value = '  +2.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.40'
$ws.Range("E9").Value = '  +12.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09008'
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.181'
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.47'
$ws.Range("E12").Value = '  -2.12%  '
$ws.Range("D13").Value = '2.103.81'
$ws.Range("E13").Value = '  -0.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.792'
$ws.Range("E14").Value = '  +0.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.817'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '96.66'
$ws.Range("E16").Value = '  -0.14%  '
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06626'
$ws.Range("E19").Value = '  -0.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.51'
$ws.Range("E20").Value = '  +1.91%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.329'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '30.531.04'
$ws.Range("E23").Value = '  -1.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.42'
$ws.Range("E24").Value = '  +0.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.353'
$ws.Range("E25").Value = '  +3.05%  '
$ws.Range("D26").Value = '2.345.87'
$ws.Range("E26").Value = '  -0.60%  '
$ws.Range("E27").Value = '  -1.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.584'
$ws.Range("E28").Value = '  -0.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '163.69'
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.98'
$ws.Range("E30").Value = '  -0.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.201'
$ws.Range("E32").Value = '  -1.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.668'
$ws.Range("E33").Value = '  +7.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.167'
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.926'
$ws.Range("E35").Value = '  -2.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.57'
$ws.Range("E36").Value = '  +10.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02583'
$ws.Range("E37").Value = '  -0.84%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.593'
$ws.Range("E38").Value = '  +0.78%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06832'
$ws.Range("E39").Value = '  +0.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2308'
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("E41").Value = '  -1.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6924'
$ws.Range("E42").Value = '  +1.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.250'
$ws.Range("E43").Value = '  +0.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.357'
$ws.Range("E44").Value = '  +5.39%  '
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '14.11'
$ws.Range("E46").Value = '  +0.31%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6395'
$ws.Range("E47").Value = '  -1.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.660'
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("E49").Value = '  -1.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.229'
$ws.Range("E50").Value = '  +3.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '83.66'
$ws.Range("E51").Value = '  +0.76%  '
